# Added transaction type 16 functionality
#
# The "RawEarnings" sheet gains a new column (TransactionType16) inserted
# immediately before the existing "ApprenticeshipContractType" column
# (which shifts one column to the right, from AB to AC). The new column
# is populated with 0 for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RawEarnings")

# Insert a new column at AB, shifting the old AB (ApprenticeshipContractType)
# and everything after it one column to the right.
$ws.Columns("AB:AB").Insert()

# Header for the newly inserted column.
$ws.Range("AB1").Value = "TransactionType16"

# Data rows (2-5) default to 0 for the new TransactionType16 column.
$ws.Range("AB2:AB5").Value = 0

# Match the column width used by the other TransactionTypeNN columns (13-27).
$ws.Columns("AB:AB").ColumnWidth = $ws.Columns("AA:AA").ColumnWidth

# Make RawEarnings the active sheet/tab, with AB1 selected - matching the
# saved view state of the edited workbook.
$ws.Activate() | Out-Null
$ws.Range("AB1").Select() | Out-Null
